# Fruta / hortaliza, semanal
# Insert a new weekly record at row 5 (pushing the existing rows 5-37 down
# to 6-38), matching the new dimension A1:R38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows down by inserting a fresh row at position 5.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(5, 1).Value  = 8
$ws.Cells.Item(5, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(5, 3).Value  = "Coquimbo"
$ws.Cells.Item(5, 4).Value  = 44687
$ws.Cells.Item(5, 5).Value  = 4
$ws.Cells.Item(5, 6).Value  = 100114007
$ws.Cells.Item(5, 7).Value  = "Jengibre"
$ws.Cells.Item(5, 8).Value  = "Sin especificar"
$ws.Cells.Item(5, 9).Value  = "Primera"
$ws.Cells.Item(5, 10).Value = 440
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14500
$ws.Cells.Item(5, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 1115
$ws.Cells.Item(5, 17).Value = 13
$ws.Cells.Item(5, 18).Value = "Hortaliza"
